# "Validacion de tolvas listo"
# Update the solver-recalculated inputs on "Cálculo por volumen" (G3, G5, G10);
# the dependent formulas on both sheets recalc automatically.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Cálculo por volumen")
$ws2 = $wb.Worksheets.Item("Cálculo de fuerza")

$ws1.Range("G3").Value = 16.4603868474641
$ws1.Range("G5").Value = 5597
$ws1.Range("G10").Value = 10.223842762403184

# Widen column H on "Cálculo por volumen" (closest reachable width to 13.7109375
# given COM's pixel-grid rounding of ColumnWidth).
$ws1.Columns.Item(8).ColumnWidth = 12.86

# Move the selection / active sheet: "Cálculo de fuerza" becomes the active
# (tab-selected) sheet with A11 selected; "Cálculo por volumen" keeps J11 selected.
$ws1.Range("J11").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("A11").Select() | Out-Null
